$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

# New localized "How to Play" description strings, added in shared-string order:
# English, German, Polish, then the key name - matching the target insertion order.
$enDesc = "The gnome runs automatically towards his house. Along the way, you can either jump or roll. There are four levels, each becoming progressively more interesting than the last.`nIn Levels 1, 2, and 3, you can jump by tapping anywhere on the screen. If you are using a computer, you can click with the mouse or use the spacebar.`nIn Level 4, you must also roll away from the owl! To jump, tap the left half of the screen, and to roll, tap the right half of the screen.`nIf you are using a computer, you can click with the mouse or use the spacebar to jump, and the Control key to roll."
$deDesc = "Der Zwerg läuft automatisch zu seinem Haus. Auf dem Weg kannst du entweder springen oder rollen. Es gibt vier Stufen, wobei jede Stufe interessanter als die vorherige wird.`nIn Stufe 1, 2 und 3 kannst du springen, indem du irgendwo auf dem Bildschirm tippst. Wenn du einen Computer verwendest, kannst du mit der Maus klicken oder die Leertaste benutzen.`nIn Stufe 4 musst du auch vor der Eule wegrollen! Um zu springen, tippe auf die linke Hälfte des Bildschirms, und um zu rollen, tippe auf die rechte Hälfte des Bildschirms.`nWenn du einen Computer verwendest, kannst du mit der Maus klicken oder die Leertaste zum Springen und die Strg-Taste zum Rollen verwenden."
$plDesc = "Krasnolud biegnie automatycznie do swojego domu. W trakcie podróży możesz skakać lub toczyć się. Są cztery poziomy, z których każdy jest coraz ciekawszy od poprzedniego.`nNa poziomach 1, 2 i 3 możesz skakać, dotykając dowolnego miejsca na ekranie. Jeśli używasz komputera, możesz kliknąć myszką lub użyć spacji.`nNa poziomie 4 musisz także toczyć się na ziemi, aby uniknąć sowy! Aby skoczyć, dotknij lewej połowy ekranu, a aby się toczyć, dotknij prawej połowy ekranu.`nJeśli używasz komputera, możesz kliknąć myszką lub użyć spacji do skakania i klawisza Ctrl do toczenia się."
$descKey = "HowToPlayDesc"

$ws.Range("B37").Value = $enDesc
$ws.Range("C37").Value = $deDesc
$ws.Range("D37").Value = $plDesc
$ws.Range("A37").Value = $descKey

# Match formatting of the row above (wrap text, same style) and row height.
$ws.Range("B37:D37").WrapText = $true
$ws.Rows.Item(37).RowHeight = 216

$ws.Range("A38").Select()
